$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of portfolio data appended below the existing table (row 42).
# Force column A to text first so the ISO-looking date string
# ("2025-09-26") is stored as literal text instead of being auto-parsed
# into a date serial number, then drop the temporary format so the cell
# keeps the sheet's default (unstyled) appearance.
$ws.Range("A42").NumberFormat = "@"
$ws.Range("A42").Value = "2025-09-26"
$ws.Range("A42").ClearFormats()

$ws.Range("B42").Value = 55.27999877929688
$ws.Range("C42").Value = 672.9000244140625
$ws.Range("D42").Value = 321
